$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 2).Value = 1398
$ws.Cells.Item(3, 3).Value = 3421
$ws.Cells.Item(3, 4).Value = 4971
$ws.Cells.Item(3, 5).Value = 28700
$ws.Cells.Item(3, 6).Value = 31600
$ws.Cells.Item(3, 7).Value = 30000

$ws.Cells.Item(8, 2).Value = 19600
$ws.Cells.Item(8, 3).Value = 30400
$ws.Cells.Item(8, 4).Value = 56200
$ws.Cells.Item(8, 5).Value = 107000
$ws.Cells.Item(8, 6).Value = 187000
$ws.Cells.Item(8, 7).Value = 286000

$ws.Cells.Item(13, 2).Value = 4751
$ws.Cells.Item(13, 3).Value = 5538
$ws.Cells.Item(13, 4).Value = 5734
$ws.Cells.Item(13, 5).Value = 5993
$ws.Cells.Item(13, 6).Value = 5514
$ws.Cells.Item(13, 7).Value = 4914

$ws.Cells.Item(18, 2).Value = 381000
$ws.Cells.Item(18, 3).Value = 484000
$ws.Cells.Item(18, 4).Value = 622000
$ws.Cells.Item(18, 5).Value = 647000
$ws.Cells.Item(18, 6).Value = 200000
$ws.Cells.Item(18, 7).Value = 122000

$ws.Cells.Item(23, 2).Value = 3312
$ws.Cells.Item(23, 3).Value = 5221
$ws.Cells.Item(23, 4).Value = 12300
$ws.Cells.Item(23, 5).Value = 7488
$ws.Cells.Item(23, 6).Value = 10500
$ws.Cells.Item(23, 7).Value = 16200

$ws.Cells.Item(28, 2).Value = 127000
$ws.Cells.Item(28, 3).Value = 161000
$ws.Cells.Item(28, 4).Value = 411000
$ws.Cells.Item(28, 5).Value = 266000
$ws.Cells.Item(28, 6).Value = 343000
$ws.Cells.Item(28, 7).Value = 463000

$ws.Cells.Item(33, 2).Value = 10800
$ws.Cells.Item(33, 3).Value = 12600
$ws.Cells.Item(33, 4).Value = 13000
$ws.Cells.Item(33, 5).Value = 15000
$ws.Cells.Item(33, 6).Value = 15200
$ws.Cells.Item(33, 7).Value = 3030

$ws.Cells.Item(38, 2).Value = 322000
$ws.Cells.Item(38, 3).Value = 386000
$ws.Cells.Item(38, 4).Value = 450000
$ws.Cells.Item(38, 5).Value = 476000
$ws.Cells.Item(38, 6).Value = 490000
$ws.Cells.Item(38, 7).Value = 492000
